# Corporate Customer excel file update:
#  - Append two new "Transaction Number" values to the bottom of the
#    single-column list on Sheet0 (rows 37 and 38).
#  - Leave the active selection on the cell the author last clicked (U6)
#    before saving, matching the workbook's saved sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New transaction numbers appended under the existing list in column A.
$ws.Range("A37").Value = "FT23185000NSGYBV"
$ws.Range("A38").Value = "FT231850JRZ2B868"

# Matches the <selection activeCell="U6" sqref="U6"/> saved in the sheet view.
$ws.Range("U6").Select()
